# Standardize "Multi KDMA" wording in the definition column (column K/AA/AD/AG/AI/AL/AM, row 2)
# Replaces the informal "multi-kdma" / "multikdma" phrasing with the standardized "Multi KDMA".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = "X refers to the block number, Y refers to the DM number the participant saw. For Eval 8 and 9, participants saw 4 blocks with 3-4 DMs each (4 in the case of Multi KDMA). The following columns describe each page of the survey using this BX_DMY format."

$ws.Range("AA2").Value = "The name and alignment value of the fourth medic being compared in this comparison page - only applies to Multi KDMA"

$ws.Range("AD2").Value = "The response to the first forced choice question (baseline vs aligned, or follow the previous column for Multi KDMA)"

$ws.Range("AG2").Value = "The difference between the Delegator|Observed_ADM alignment comparison between the aligned ADM and the baseline ADM (or whatever is in the Alignment column for Multi KDMA)"

$ws.Range("AI2").Value = "The response to the second forced choice question (aligned vs misaligned, or follow the previous column for Multi KDMA)"

$ws.Range("AL2").Value = "The difference between the Delegator|Observed_ADM alignment comparison between the aligned ADM and the misaligned ADM (or follow the alignment column for Multi KDMA)"

$ws.Range("AM2").Value = "The alignment of the third two DMs being compared (Multi KDMA only)"
